$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the style from the existing H1 header
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-38
$iValues = @(8,6,8,9,7,7,5,8,4,8,6,6,7,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,5)
$jValues = @(8,6,9,9,7,8,6,9,6,8,6,7,7,3,6,6,7,4,7,6,7,4,7,2,7,6,4,7,7,6,7,7,5,7,6,5,7)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}

$ws.Range("A1").Select() | Out-Null
